$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E9").Value = 13.72050000000001
$ws.Range("E13").Value = 12.3824
$ws.Range("E16").Value = 12.52370000000001
$ws.Range("E18").Value = 12.8457
$ws.Range("E20").Value = 13.13939999999999
